# Update crypto price (D) and 1h volume-change (E) columns per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.068.45"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "2.472.17"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'577.67"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "'146.77"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("E8").Value = "  -0.64%  "
$ws.Range("D9").Value = "2.471.05"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("D10").Value = "'0.112"
$ws.Range("E10").Value = "  +0.57%  "
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").Value = "'0.353"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").Value = "'28.74"
$ws.Range("E14").Value = "  +5.06%  "
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").Value = "2.919.79"
$ws.Range("E16").Value = "  +2.16%  "
$ws.Range("D17").Value = "62.943.54"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "2.477.68"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").Value = "'8.21"
$ws.Range("E19").Value = "  +3.88%  "
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("D21").Value = "'329.53"
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("D22").Value = "'2.24"
$ws.Range("E22").Value = "  +9.66%  "
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").Value = "'66.34"
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("D26").Value = "'668.50"
$ws.Range("E26").Value = "  +5.06%  "
$ws.Range("D27").Value = "'9.62"
$ws.Range("E27").Value = "  +12.19%  "
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").Value = "2.593.17"
$ws.Range("E29").Value = "  +2.18%  "
$ws.Range("E30").Value = "  -9.30%  "
$ws.Range("E31").Value = "  +2.81%  "
$ws.Range("D32").Value = "'8.06"
$ws.Range("E32").Value = "  -1.66%  "
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("E34").Value = "  -3.67%  "
$ws.Range("E35").Value = "  +3.37%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("D38").Value = "'5.48"
$ws.Range("E38").Value = "  +1.51%  "
$ws.Range("D39").Value = "'0.371"
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("D40").Value = "'18.78"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("D41").Value = "'150.46"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("E45").Value = "  +6.04%  "
$ws.Range("D46").Value = "'152.23"
$ws.Range("E46").Value = "  +5.07%  "
$ws.Range("E47").Value = "  +19.37%  "
$ws.Range("D48").Value = "'3.60"
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").Value = "'20.64"
$ws.Range("E49").Value = "  +0.99%  "
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("E51").Value = "  -0.80%  "
